# Repopulate the BirdFamilies "Sheet1" table from AvesTables.xlsx.
# The old ad-hoc rows (many with typos / outdated taxonomy, e.g. "Caprimuldidae",
# "Scolopacide") are replaced wholesale with the authoritative 88-row family list
# used elsewhere, keyed the same way: column A = Id, column B = Name, column C = SciName.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row stays the same, but make sure it is exactly right.
$ws.Cells.Item(1, 1).Value = "Id"
$ws.Cells.Item(1, 2).Value = "Name"
$ws.Cells.Item(1, 3).Value = "SciName"

# Full replacement data set (Id, Name, SciName), 88 bird families.
$data = @(
    @('1', 'Hawks and Eagles', 'Accipitridae'),
    @('2', 'Long-tailed Tits', 'Aegithalidae'),
    @('3', 'Larks', 'Alaudidae'),
    @('4', 'Kingfishers and Allies', 'Alcedinidae'),
    @('5', 'Auks, Murres and Puffins', 'Alcidae'),
    @('6', 'Ducks, Geese and Swans', 'Anatidae'),
    @('7', 'Anhingas', 'Anhingidae'),
    @('8', 'Swifts', 'Apodidae'),
    @('9', 'Limpkin', 'Aramidae'),
    @('10', 'Herons, Egrets and Bitterns', 'Ardeidae'),
    @('11', 'Waxwings', 'Bombycillidae'),
    @('12', 'Thick-knees', 'Burhinidae'),
    @('13', 'Longspurs and Snow Buntings', 'Calcariidae'),
    @('14', 'Nightjars and Allies', 'Caprimulgidae'),
    @('15', 'Cardinals and Allies', 'Cardinalidae'),
    @('16', 'New World Vultures', 'Cathartidae'),
    @('17', 'Treekeepers', 'Certhiidae'),
    @('18', 'Plovers and Lapwings', 'Charadriidae'),
    @('19', 'Storks', 'Ciconiidae'),
    @('20', 'Dippers', 'Cinclidae'),
    @('21', 'Pigeons and Doves', 'Columbidae'),
    @('22', 'Crows, Jays and Magpies', 'Corvidae'),
    @('23', 'Guans, Chachalacas and Curassows', 'Cracidae'),
    @('24', 'Cuckoos', 'Cuculidae'),
    @('25', 'Albatrosses', 'Diomedeidae'),
    @('26', 'Buntings', 'Emberizidae'),
    @('27', 'Falcons and Caracaras', 'Falconidae'),
    @('28', 'Frigatebirds', 'Fregatidae'),
    @('29', 'Finches, Euphorias and Allies', 'Fringillidae'),
    @('30', 'Loons', 'Gaviidae'),
    @('31', 'Coursers and Pratincoles', 'Glareolidae'),
    @('32', 'Cranes', 'Gruidae'),
    @('33', 'Oystercatchers', 'Haematopodidae'),
    @('34', 'Swallows and Martins', 'Hirundinidae'),
    @('35', 'Northern Storm-Petrels', 'Hydrobatidae'),
    @('36', 'Troupials and Allies', 'Icteridae'),
    @('37', 'Yellow-breasted Chat', 'Icteriidae'),
    @('38', 'Shrikes', 'Laniidae'),
    @('39', 'Gulls, Terns and Skimmers', 'Laridae'),
    @('40', 'Mockingbirds and Thrashers', 'Mimidae'),
    @('41', 'Wagtails and Pipits', 'Motacillidae'),
    @('42', 'Old World Flycatchers', 'Muscicapidae'),
    @('43', 'New World Quail', 'Odontophoridae'),
    @('44', 'Orioles', 'Oriolidae'),
    @('45', 'Bustards', 'Otididae'),
    @('46', 'Osprey', 'Pandionidae'),
    @('47', 'Tits, Chickadees and Titmice', 'Paridae'),
    @('48', 'New World Warblers', 'Parulidae'),
    @('49', 'New World Sparrows', 'Passerellidae'),
    @('50', 'Old World Sparrows', 'Passeridae'),
    @('51', 'Pelicans', 'Pelecanidae'),
    @('52', 'Olive Warbler', 'Peucedramidae'),
    @('53', 'TropicBirds', 'Phaethontidae'),
    @('54', 'Cormorants and Shags', 'Phalacrocoridae'),
    @('55', 'Phalaropes', 'Phalaropidae'),
    @('56', 'Pheasants and Partridges', 'Phasianidae'),
    @('57', 'Flamingos', 'Phoenicopteridae'),
    @('58', 'Leaf Warblers', 'Phylloscopidae'),
    @('59', 'Woodpeckers', 'Picidae'),
    @('60', 'Grebes', 'Podicipedidae'),
    @('61', 'Gnatcatchers', 'Polioptilidae'),
    @('62', 'Shearwaters and Petrels', 'Procellariidae'),
    @('63', 'Accentors', 'Prunellidae'),
    @('64', 'Parrots', 'Psittacidae'),
    @('65', 'Sandgrouse', 'Pteroclididae'),
    @('66', 'Silky-flycatchers', 'Ptiliogonatidae'),
    @('67', 'Rails, Gallinules and Coots', 'Rallidae'),
    @('68', 'Avocets ans Stilts', 'Recurvirostridae'),
    @('69', 'Kinglets', 'Regulidae'),
    @('70', 'Penduline-Tits', 'Remizidae'),
    @('71', 'Sandpipers ans Allies', 'Scolopacidae'),
    @('72', 'Nuthatches', 'Sittidae'),
    @('73', 'Skuas and Jaegers', 'Stercorariidae'),
    @('74', 'Terns', 'Sternidae'),
    @('75', 'Owls', 'Strigidae'),
    @('76', 'Starlings', 'Sturnidae'),
    @('77', 'Gannets and Boobies', 'Sulidae'),
    @('78', 'Sylviid Warblers, Parrotbills and Allies', 'Sylviidae'),
    @('79', 'Grouse', 'Tetraonidae'),
    @('80', 'Tanagers', 'Thraupidae'),
    @('81', 'Ibises and Spoonbills', 'Threskiornithidae'),
    @('82', 'Hummingbirds', 'Trochilidae'),
    @('83', 'Wrens', 'Troglodytidae'),
    @('84', 'Trogons', 'Trogonidae'),
    @('85', 'Thrushes and Allies', 'Turdidae'),
    @('86', 'Tyrant Flycatchers', 'Tyrannidae'),
    @('87', 'Barn-Owls', 'Tytonidae'),
    @('88', 'Vireos, Shrike-Babblers and Erpornis', 'Vireonidae')
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = [int]$row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Reset the view: select A1 (previously the sheet was scrolled down with A68 selected).
$ws.Range("A1").Select()
